# Auto-generated edit script: refresh computed market-price columns (H-N)
# across multiple job sheets, per the scheduled-runner data update.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 797.3125
$ws.Range("I19").Value = 776.4286
$ws.Range("J19").Value = 813.55554
$ws.Range("K19").Value = 776.4286
$ws.Range("L19").Value = 813.55554
$ws.Range("M19").Value = -601.4286
$ws.Range("N19").Value = -1163.55554
$ws.Range("H103").Value = 913.875
$ws.Range("I103").Value = 519.5
$ws.Range("J103").Value = 934.6316
$ws.Range("K103").Value = 1558.5
$ws.Range("L103").Value = 2803.8948
$ws.Range("M103").Value = -972.5
$ws.Range("N103").Value = -3975.8948
$ws.Range("H107").Value = 396.22223
$ws.Range("I107").Value = 382.18182
$ws.Range("J107").Value = 458
$ws.Range("K107").Value = 382.18182
$ws.Range("L107").Value = 458
$ws.Range("M107").Value = 1537.81818
$ws.Range("N107").Value = -4298
$ws.Range("H129").Value = 4379.759
$ws.Range("J129").Value = 963.44446
$ws.Range("L129").Value = 2890.33338
$ws.Range("N129").Value = -12890.33338

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 34597.71
$ws.Range("I32").Value = 5296.977
$ws.Range("K32").Value = 5296.977
$ws.Range("M32").Value = -5009.977
$ws.Range("H45").Value = 2332.7058
$ws.Range("I45").Value = 1695.2
$ws.Range("J45").Value = 3243.4285
$ws.Range("K45").Value = 1695.2
$ws.Range("L45").Value = 3243.4285
$ws.Range("M45").Value = -1318.2
$ws.Range("N45").Value = -3997.4285
$ws.Range("H88").Value = 3173.875
$ws.Range("I88").Value = 3231.8333
$ws.Range("J88").Value = 3000
$ws.Range("K88").Value = 3231.8333
$ws.Range("L88").Value = 3000
$ws.Range("M88").Value = -2825.8333
$ws.Range("N88").Value = -3812
$ws.Range("H91").Value = 3173.875
$ws.Range("I91").Value = 3231.8333
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 3231.8333
$ws.Range("L91").Value = 3000
$ws.Range("M91").Value = -1827.8333
$ws.Range("N91").Value = -5808
$ws.Range("H102").Value = 37441.965
$ws.Range("I102").Value = 60394
$ws.Range("J102").Value = 1970.6364
$ws.Range("K102").Value = 60394
$ws.Range("L102").Value = 1970.6364
$ws.Range("M102").Value = -58772
$ws.Range("N102").Value = -5214.6364

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1510.5625
$ws.Range("I99").Value = 1262.7273
$ws.Range("J99").Value = 2055.8
$ws.Range("K99").Value = 1262.7273
$ws.Range("L99").Value = 2055.8
$ws.Range("M99").Value = 235.2727
$ws.Range("N99").Value = -5051.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 654.5454999999999
$ws.Range("I16").Value = 318.5
$ws.Range("J16").Value = 1057.8
$ws.Range("K16").Value = 318.5
$ws.Range("L16").Value = 1057.8
$ws.Range("M16").Value = -31.5
$ws.Range("N16").Value = -1631.8
$ws.Range("H31").Value = 22758.307
$ws.Range("I31").Value = 43047.625
$ws.Range("J31").Value = 3280.56
$ws.Range("K31").Value = 43047.625
$ws.Range("L31").Value = 3280.56
$ws.Range("M31").Value = -42752.625
$ws.Range("N31").Value = -3870.56
$ws.Range("H34").Value = 22758.307
$ws.Range("I34").Value = 43047.625
$ws.Range("J34").Value = 3280.56
$ws.Range("K34").Value = 43047.625
$ws.Range("L34").Value = 3280.56
$ws.Range("M34").Value = -42845.625
$ws.Range("N34").Value = -3684.56
$ws.Range("H62").Value = 2571.4285
$ws.Range("J62").Value = 2680
$ws.Range("L62").Value = 2680
$ws.Range("N62").Value = -3928
$ws.Range("H65").Value = 2571.4285
$ws.Range("J65").Value = 2680
$ws.Range("L65").Value = 13400
$ws.Range("N65").Value = -19640
$ws.Range("H113").Value = 654.5454999999999
$ws.Range("I113").Value = 318.5
$ws.Range("J113").Value = 1057.8
$ws.Range("K113").Value = 318.5
$ws.Range("L113").Value = 1057.8
$ws.Range("M113").Value = 1851.5
$ws.Range("N113").Value = -5397.8
$ws.Range("H122").Value = 492.14285
$ws.Range("I122").Value = 435.45456
$ws.Range("J122").Value = 700
$ws.Range("K122").Value = 1306.36368
$ws.Range("L122").Value = 2100
$ws.Range("M122").Value = 1143.63632
$ws.Range("N122").Value = -7000
$ws.Range("H132").Value = 36588748
$ws.Range("I132").Value = 33336442
$ws.Range("J132").Value = 45458668
$ws.Range("K132").Value = 100009326
$ws.Range("L132").Value = 136376004
$ws.Range("M132").Value = -100006796
$ws.Range("N132").Value = -136381064
$ws.Range("H134").Value = 1755.3182
$ws.Range("I134").Value = 1947.4615
$ws.Range("J134").Value = 1477.7778
$ws.Range("K134").Value = 5842.3845
$ws.Range("L134").Value = 4433.3334
$ws.Range("M134").Value = -3307.3845
$ws.Range("N134").Value = -9503.3334

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 612.85297
$ws.Range("I113").Value = 567.1111
$ws.Range("J113").Value = 664.3125
$ws.Range("K113").Value = 1701.3333
$ws.Range("L113").Value = 1992.9375
$ws.Range("M113").Value = 468.6667000000002
$ws.Range("N113").Value = -6332.9375
$ws.Range("H131").Value = 816.58
$ws.Range("I131").Value = 490.76923
$ws.Range("J131").Value = 865.2643399999999
$ws.Range("K131").Value = 1472.30769
$ws.Range("L131").Value = 2595.79302
$ws.Range("M131").Value = 3567.69231
$ws.Range("N131").Value = -12675.79302
$ws.Range("H132").Value = 2975.6
$ws.Range("J132").Value = 2559.0833
$ws.Range("L132").Value = 23031.7497
$ws.Range("N132").Value = -28091.7497

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 13996.667
$ws.Range("J62").Value = 14996
$ws.Range("L62").Value = 14996
$ws.Range("N62").Value = -16368
$ws.Range("H65").Value = 13996.667
$ws.Range("J65").Value = 14996
$ws.Range("L65").Value = 44988
$ws.Range("N65").Value = -51852

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 126463.125
$ws.Range("I16").Value = 167283.5
$ws.Range("K16").Value = 167283.5
$ws.Range("M16").Value = -167113.5
$ws.Range("H22").Value = 968.8125
$ws.Range("I22").Value = 1149.5
$ws.Range("J22").Value = 943
$ws.Range("K22").Value = 1149.5
$ws.Range("L22").Value = 943
$ws.Range("M22").Value = -854.5
$ws.Range("N22").Value = -1533
$ws.Range("H27").Value = 968.8125
$ws.Range("I27").Value = 1149.5
$ws.Range("J27").Value = 943
$ws.Range("K27").Value = 1149.5
$ws.Range("L27").Value = 943
$ws.Range("M27").Value = -1042.5
$ws.Range("N27").Value = -1157
$ws.Range("H46").Value = 3149.889
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 3149.889
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 3149.889
$ws.Range("N46").Value = -3525.889
$ws.Range("M46").Value = $null
$ws.Range("H61").Value = 2516.6924
$ws.Range("I61").Value = 1460
$ws.Range("K61").Value = 1460
$ws.Range("M61").Value = -1258
$ws.Range("H82").Value = 1932
$ws.Range("I82").Value = 2300
$ws.Range("J82").Value = 1441.3334
$ws.Range("K82").Value = 2300
$ws.Range("L82").Value = 1441.3334
$ws.Range("M82").Value = -1939
$ws.Range("N82").Value = -2163.3334
$ws.Range("H85").Value = 1932
$ws.Range("I85").Value = 2300
$ws.Range("J85").Value = 1441.3334
$ws.Range("K85").Value = 2300
$ws.Range("L85").Value = 1441.3334
$ws.Range("M85").Value = -1052
$ws.Range("N85").Value = -3937.3334
$ws.Range("H113").Value = 2516.6924
$ws.Range("I113").Value = 1460
$ws.Range("K113").Value = 1460
$ws.Range("M113").Value = 710
$ws.Range("H132").Value = 1754.3914
$ws.Range("I132").Value = 1496.5555
$ws.Range("J132").Value = 2682.6
$ws.Range("K132").Value = 4489.666499999999
$ws.Range("L132").Value = 8047.799999999999
$ws.Range("M132").Value = -1959.666499999999
$ws.Range("N132").Value = -13107.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 33000
$ws.Range("J69").Value = 33000
$ws.Range("L69").Value = 33000
$ws.Range("N69").Value = -34498
$ws.Range("H72").Value = 33000
$ws.Range("J72").Value = 33000
$ws.Range("L72").Value = 99000
$ws.Range("N72").Value = -106488
$ws.Range("H96").Value = 200001630
$ws.Range("I96").Value = 333335000
$ws.Range("J96").Value = 1599.5
$ws.Range("K96").Value = 333335000
$ws.Range("L96").Value = 1599.5
$ws.Range("M96").Value = -333333627
$ws.Range("N96").Value = -4345.5
$ws.Range("H122").Value = 1810.1
$ws.Range("I122").Value = 946.4666999999999
$ws.Range("J122").Value = 4401
$ws.Range("K122").Value = 2839.4001
$ws.Range("L122").Value = 13203
$ws.Range("M122").Value = -389.4000999999998
$ws.Range("N122").Value = -18103

